$d = $word.ActiveDocument
$p = $d.Paragraphs(22)
try { $p.Range.ListFormat.ListType } catch { Write-Output "ListType err: $_" }
try { $p.Style = "Normal" ; Write-Output "style set ok"} catch { Write-Output "style err: $_" }
